# Updated cryptos list on Mon Sep 11 03:15:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.874.43'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '1.620.44'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.95%  '
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '1.845.15'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").Value = '1.628.92'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("E14").Value = '  -2.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").Value = '25.874.37'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  -2.82%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("E27").Value = '  -4.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("E30").Value = '  -1.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0478'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.29%  '
$ws.Range("E34").Value = '  -2.10%  '
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("D36").Value = '1.124.21'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.839'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.62%  '
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0153'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.511'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").Value = '1.755.54'
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.751'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.65%  '
$ws.Range("E44").Value = '  -5.37%  '
$ws.Range("D45").Value = '0.0₆0112'
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.76%  '